# Automatische test-sync: 2025-06-22 18:54:50
# Adds the newly-logged inbound mail ("Offerte voor 500 stuks") as row 24 on
# the "Logs" sheet, extends the conditional formatting ranges to cover it,
# and refreshes the "Dashboard" category-count table (which is sorted by
# descending count) now that "Offerte / Prijsaanvraag" has climbed from 1 to 2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new mail-log entry as row 24
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A24").Value = "Offerte voor 500 stuks"
$logs.Range("B24").Value = "mailmind.test@zohomail.eu"
$logs.Range("C24").Value = "Graag ontvang ik een offerte voor 500 stuks van product X."
$logs.Range("D24").Value = "Offerte / Prijsaanvraag"
$logs.Range("E24").Value = "Beste klant,`nDank u wel voor uw interesse in product X. Om u een offerte op maat te kunnen sturen, hebben wij wat meer informatie nodig. Kunt u ons laten weten welke specifieke variant(en) van product X u wenst te bestellen en in welke regio de levering zal plaatsvinden? Op basis van deze gegevens kunnen wij een passende offerte voor u opstellen.`nAlvast bedankt voor de aanvullende informatie.`nMet vriendelijke groet,`n[Naam Bedrijf]"
$logs.Range("F24").Value = "2025-06-22 18:54:34"
$logs.Range("G24").Value = "Ja"

# ---------------------------------------------------------------------
# 2. Logs sheet: widen the conditional-formatting ranges so the new row
#    (D24/G24) is covered, same as the D2:D23 -> D2:D24 / G2:G23 -> G2:G24
#    range bump Excel performs when a table-backed range grows.
# ---------------------------------------------------------------------
$catRules = $logs.Range("D2:D23").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D24"))
}

$answeredRules = $logs.Range("G2:G23").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G24"))
}

# ---------------------------------------------------------------------
# 3. Dashboard sheet: "Offerte / Prijsaanvraag" now has 2 hits instead of
#    1, so the count table (sorted descending by Aantal) is re-ordered.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Offerte / Prijsaanvraag"
$dash.Range("B5").Value = 2

$dash.Range("A6").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B6").Value = 2

$dash.Range("A7").Value = "Overig"
$dash.Range("B7").Value = 2

$dash.Range("A8").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B8").Value = 2

$dash.Range("A9").Value = "Klacht / Probleem"
$dash.Range("B9").Value = 1

$dash.Range("A10").Value = "Uitnodiging / Evenement"
$dash.Range("B10").Value = 1

$dash.Range("A11").Value = "Openingstijden / Locatie"
$dash.Range("B11").Value = 1

$dash.Range("A12").Value = "Retour / Terugbetaling"
$dash.Range("B12").Value = 1
